# Update cryptos list prices/volumes (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '64.056.48'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.31%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.758.18'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  +0.06%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '577.51'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.66%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '158.82'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.66%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.26%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.605'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -2.81%  '
$ws.Range("E9").Value = '  -2.01%  '
$ws.Range("E10").Value = '  +1.87%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.384'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -2.53%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '5.63'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -16.91%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '3.242.60'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.08%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '26.94'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.91%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '63.696.41'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("E16").Value = '  -1.79%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.755.10'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.09%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '12.20'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.88%  '
$ws.Range("E19").Value = '  -1.86%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '357.78'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.50%  '
$ws.Range("E21").Value = '  -3.24%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.02%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.536'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.48%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '65.42'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.19%  '
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  +0.27%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.0₃0910'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.98%  '
$ws.Range("E29").Value = '  -3.23%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '7.25'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.22%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.26'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.23%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '169.85'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.95%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.97'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.87%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '20.25'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.97%  '
$ws.Range("E35").Value = '  +1.03%  '
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("E37").Value = '  -1.13%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.46%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '6.30'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.14%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '337.37'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.50%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '4.20'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -2.58%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '39.15'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.85%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '21.52'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.02%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '21.83'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.12%  '
$ws.Range("E45").Value = '  -2.62%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0257'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.74%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.103'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.15%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '135.56'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.16%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.628'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -3.27%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("E51").Value = '  +0.04%  '
